# eCoaching pilot-question update
# TFS 9511 - New pilot question in CSR Survey
#
# Adds a 6th ("pilot") survey question end-to-end:
#  - Survey_DIM_Question: new isPilot column + 2 new question rows
#  - Survey_DIM_Response: 5 new response rows
#  - Survey_DIM_QAnswer:  new isPilot/ResponseOrder columns + 7 new rows
#  - New "Survey_Sites" worksheet (mirrors DIM_Site, adds isPilot/isHotTopic flags)
#  - Revision_History: new row documenting the change

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Survey_DIM_Question (sheet18) - add isPilot column + 2 new rows
# ------------------------------------------------------------------
$wsQ = $wb.Worksheets.Item("Survey_DIM_Question")

# New column I ("isPilot") - highlighted (yellow fill) header + existing rows
$wsQ.Range("I1").Value = "isPilot"
$wsQ.Range("I1:I7").Interior.Color = 65535
$wsQ.Range("I2").Value = 0
$wsQ.Range("I3").Value = 0
$wsQ.Range("I4").Value = 0
$wsQ.Range("I5").Value = 0
$wsQ.Range("I6").Value = 0
$wsQ.Range("I7").Value = 0

# New row 9 first (pilot question) so the shared-string table receives
# "How prepared..." (903) before "Hot Topic question" (904)
$wsQ.Range("A9").Value = 7
$wsQ.Range("B9").Value = "How prepared was your supervisor during your coaching session?| Please explain below."
$wsQ.Range("C9").Value = 7
$wsQ.Range("D9").Value = 20180201
$wsQ.Range("E9").Value = 99991231
$wsQ.Range("F9").Value = 0
$wsQ.Range("G9").Value = 1
$wsQ.Range("H7").Copy($wsQ.Range("H9"))
$wsQ.Range("H9").Value = 43132
$wsQ.Range("H9").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss:mss"
$wsQ.Range("H9").Interior.Color = 65535
$wsQ.Range("I9").Value = 1
$wsQ.Range("A9:G9").Interior.Color = 65535
$wsQ.Range("I9").Interior.Color = 65535

# New row 8 (hot topic question - not pilot, not highlighted)
$wsQ.Range("A8").Value = 6
$wsQ.Range("B8").Value = "Hot Topic question"
$wsQ.Range("C8").Value = 6
$wsQ.Range("D8").Value = 20150901
$wsQ.Range("E8").Value = 20150930
$wsQ.Range("F8").Value = 1
$wsQ.Range("G8").Value = 0
$wsQ.Range("H7").Copy($wsQ.Range("H8"))
$wsQ.Range("H8").Value = 42248
$wsQ.Range("I8").Value = 0
$wsQ.Range("I8").Interior.Color = 65535

$wsQ.Range("I1:I9").Select()
